# docs/counter.xlsx — "Fixed form to use /data root node (#21)"
#
# The ODK "counter2" question's appearance column referenced the survey's
# root node via the literal sheet/form name ("/counter/...") instead of
# the canonical XLSForm root "/data/...". Correct the itext() xpath
# references in the shared text so the generated XForm resolves them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two broken root-node references inside the "appearance" cell
# (row 3, column F) of the survey sheet:
#   jr:itext('/counter/form_name:label') -> jr:itext('/data/form_name:label')
#   jr:itext('/counter/counter2:label')  -> jr:itext('/data/counter2:label')
$null = $ws.Cells.Replace("/counter/form_name:label", "/data/form_name:label")
$null = $ws.Cells.Replace("/counter/counter2:label", "/data/counter2:label")

# Restore the cursor/selection state captured in the saved workbook view.
$null = $ws.Range("F22").Select()

# Best-effort view cosmetics captured by the original author's save
# (tab-bar split ratio and the window's top-left scroll anchor). Not all
# hosts persist these window-chrome settings, but set them anyway so the
# intent is preserved if the runtime supports it.
try { $excel.ActiveWindow.TabRatio = 500 } catch {}
try {
    $excel.ActiveWindow.ScrollColumn = 5
    $excel.ActiveWindow.ScrollRow = 1
} catch {}
